# Capitalize the hex-digit letters (a-f -> A-F) inside the "doip" (G) and
# "uds" (H) columns' byte-code strings, e.g. "0x02:0xfd:...:0x0e" becomes
# "0x02:0xFD:...:0x0E". The "0x" prefix itself is left as lowercase "x";
# only the two hex digits following it are affected. Cells such as "N/A"
# are unaffected since they contain no a-f letters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 41

for ($row = 2; $row -le $lastRow; $row++) {
    foreach ($col in @("G", "H")) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        $val = $cell.Value()
        if ($val -ne $null) {
            $newVal = $val.Replace("a", "A").Replace("b", "B").Replace("c", "C").Replace("d", "D").Replace("e", "E").Replace("f", "F")
            # Note: avoid comparing $newVal/$val with -eq/-ne here, since string
            # comparisons in this runtime are case-insensitive; just assign
            # unconditionally (values such as "N/A" are unaffected by Replace).
            $cell.Value = $newVal
        }
    }
}

Write-Output "Capitalized hex letters in columns G and H for rows 2-$lastRow"
